# CompStat weekly data refresh: new crime data collected.
# Update report header (volume number + week-covering date range).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# Stable anchor cells (outside the edited rows 14-30) used to copy the
# correct cell style+type when a cell's category flips between the
# "no data" text placeholder ("0" / "***.*") and a real number.
$anchorZeroText = $ws.Range("C33")     # style 13, text "0"   (shared string idx 20)
$anchorStarText = $ws.Range("E33")     # style 13, text "***.*" (shared string idx 21)
$anchorCountNum = $ws.Range("I31")     # style 14, numeric count format
$anchorPctNum   = $ws.Range("L31")     # style 15, numeric % change format


$anchorCountNum.Copy($ws.Cells.Item(14,4))
$ws.Cells.Item(14,4).Value = 1
$anchorPctNum.Copy($ws.Cells.Item(14,5))
$ws.Cells.Item(14,5).Value = -100
$anchorCountNum.Copy($ws.Cells.Item(14,7))
$ws.Cells.Item(14,7).Value = 1
$anchorPctNum.Copy($ws.Cells.Item(14,8))
$ws.Cells.Item(14,8).Value = -100
$ws.Cells.Item(14,10).Value = 2
$ws.Cells.Item(14,11).Value = -50
$ws.Cells.Item(14,14).Value = -87.5
$anchorZeroText.Copy($ws.Cells.Item(15,4))
$anchorStarText.Copy($ws.Cells.Item(15,5))
$anchorZeroText.Copy($ws.Cells.Item(15,6))
$ws.Cells.Item(15,8).Value = -100
$ws.Cells.Item(15,14).Value = -40
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 3
$ws.Cells.Item(16,5).Value = 66.666666666666
$ws.Cells.Item(16,6).Value = 9
$ws.Cells.Item(16,7).Value = 15
$ws.Cells.Item(16,8).Value = -40
$ws.Cells.Item(16,9).Value = 22
$ws.Cells.Item(16,10).Value = 31
$ws.Cells.Item(16,11).Value = -29.032258064516
$ws.Cells.Item(16,12).Value = -18.518518518518
$ws.Cells.Item(16,13).Value = -46.341463414634
$ws.Cells.Item(16,14).Value = -84.057971014492
$ws.Cells.Item(17,3).Value = 7
$ws.Cells.Item(17,4).Value = 6
$ws.Cells.Item(17,5).Value = 16.666666666666
$ws.Cells.Item(17,6).Value = 17
$ws.Cells.Item(17,7).Value = 17
$ws.Cells.Item(17,8).Value = 0
$ws.Cells.Item(17,9).Value = 31
$ws.Cells.Item(17,10).Value = 38
$ws.Cells.Item(17,11).Value = -18.421052631578
$ws.Cells.Item(17,12).Value = -26.190476190476
$ws.Cells.Item(17,13).Value = -26.190476190476
$ws.Cells.Item(17,14).Value = -67.708333333333
$ws.Cells.Item(18,4).Value = 6
$ws.Cells.Item(18,5).Value = -50
$ws.Cells.Item(18,6).Value = 7
$ws.Cells.Item(18,7).Value = 10
$ws.Cells.Item(18,8).Value = -30
$ws.Cells.Item(18,9).Value = 14
$ws.Cells.Item(18,10).Value = 19
$ws.Cells.Item(18,11).Value = -26.315789473684
$ws.Cells.Item(18,12).Value = -6.666666666666
$ws.Cells.Item(18,13).Value = -22.222222222222
$ws.Cells.Item(18,14).Value = -89.629629629629
$ws.Cells.Item(19,3).Value = 3
$ws.Cells.Item(19,4).Value = 13
$ws.Cells.Item(19,5).Value = -76.923076923076
$ws.Cells.Item(19,7).Value = 35
$ws.Cells.Item(19,8).Value = -37.142857142857
$ws.Cells.Item(19,9).Value = 39
$ws.Cells.Item(19,10).Value = 79
$ws.Cells.Item(19,11).Value = -50.632911392405
$ws.Cells.Item(19,12).Value = -38.095238095238
$ws.Cells.Item(19,13).Value = -18.75
$ws.Cells.Item(19,14).Value = -45.833333333333
$anchorZeroText.Copy($ws.Cells.Item(20,3))
$anchorZeroText.Copy($ws.Cells.Item(20,4))
$anchorStarText.Copy($ws.Cells.Item(20,5))
$ws.Cells.Item(20,7).Value = 5
$ws.Cells.Item(20,8).Value = -60
$ws.Cells.Item(20,12).Value = -37.5
$ws.Cells.Item(20,14).Value = -81.481481481481
$ws.Cells.Item(21,4).Value = 29
$ws.Cells.Item(21,5).Value = -37.931034482758
$ws.Cells.Item(21,6).Value = 57
$ws.Cells.Item(21,7).Value = 84
$ws.Cells.Item(21,8).Value = -32.142857142857
$ws.Cells.Item(21,9).Value = 115
$ws.Cells.Item(21,10).Value = 180
$ws.Cells.Item(21,11).Value = -36.111111111111
$ws.Cells.Item(21,12).Value = -27.672955974842
$ws.Cells.Item(21,13).Value = -24.342105263157
$ws.Cells.Item(21,14).Value = -76.091476091476
$ws.Cells.Item(22,13).Value = -50
$ws.Cells.Item(23,3).Value = 1
$anchorCountNum.Copy($ws.Cells.Item(23,4))
$ws.Cells.Item(23,4).Value = 3
$anchorPctNum.Copy($ws.Cells.Item(23,5))
$ws.Cells.Item(23,5).Value = -66.666666666666
$ws.Cells.Item(23,6).Value = 5
$ws.Cells.Item(23,7).Value = 6
$ws.Cells.Item(23,8).Value = -16.666666666666
$ws.Cells.Item(23,9).Value = 15
$ws.Cells.Item(23,10).Value = 18
$ws.Cells.Item(23,11).Value = -16.666666666666
$ws.Cells.Item(23,12).Value = 25
$ws.Cells.Item(23,13).Value = 66.666666666666
$ws.Cells.Item(24,3).Value = 26
$ws.Cells.Item(24,4).Value = 22
$ws.Cells.Item(24,5).Value = 18.181818181818
$ws.Cells.Item(24,6).Value = 142
$ws.Cells.Item(24,8).Value = 67.058823529411
$ws.Cells.Item(24,9).Value = 280
$ws.Cells.Item(24,10).Value = 163
$ws.Cells.Item(24,11).Value = 71.779141104294
$ws.Cells.Item(24,12).Value = 40
$ws.Cells.Item(24,13).Value = 105.882352941176
$ws.Cells.Item(25,3).Value = 17
$ws.Cells.Item(25,5).Value = 30.769230769230
$ws.Cells.Item(25,6).Value = 103
$ws.Cells.Item(25,7).Value = 53
$ws.Cells.Item(25,8).Value = 94.339622641509
$ws.Cells.Item(25,9).Value = 198
$ws.Cells.Item(25,10).Value = 84
$ws.Cells.Item(25,11).Value = 135.714285714286
$ws.Cells.Item(25,12).Value = 37.5
$ws.Cells.Item(26,3).Value = 8
$ws.Cells.Item(26,4).Value = 5
$ws.Cells.Item(26,5).Value = 60
$ws.Cells.Item(26,6).Value = 34
$ws.Cells.Item(26,7).Value = 32
$ws.Cells.Item(26,8).Value = 6.25
$ws.Cells.Item(26,9).Value = 75
$ws.Cells.Item(26,10).Value = 66
$ws.Cells.Item(26,11).Value = 13.636363636363
$ws.Cells.Item(26,12).Value = 25
$ws.Cells.Item(26,13).Value = 5.633802816901
$anchorZeroText.Copy($ws.Cells.Item(27,4))
$anchorStarText.Copy($ws.Cells.Item(27,5))
$anchorZeroText.Copy($ws.Cells.Item(27,6))
$ws.Cells.Item(27,8).Value = -100
$ws.Cells.Item(27,12).Value = 0
$anchorCountNum.Copy($ws.Cells.Item(28,3))
$ws.Cells.Item(28,3).Value = 1
$ws.Cells.Item(28,4).Value = 2
$ws.Cells.Item(28,5).Value = -50
$ws.Cells.Item(28,6).Value = 3
$ws.Cells.Item(28,7).Value = 4
$ws.Cells.Item(28,8).Value = -25
$ws.Cells.Item(28,9).Value = 5
$ws.Cells.Item(28,10).Value = 7
$ws.Cells.Item(28,11).Value = -28.571428571428
$ws.Cells.Item(28,12).Value = -28.571428571428
$anchorCountNum.Copy($ws.Cells.Item(29,4))
$ws.Cells.Item(29,4).Value = 1
$anchorPctNum.Copy($ws.Cells.Item(29,5))
$ws.Cells.Item(29,5).Value = -100
$ws.Cells.Item(29,7).Value = 2
$ws.Cells.Item(29,10).Value = 2
$ws.Cells.Item(29,11).Value = -50
$ws.Cells.Item(29,14).Value = -94.736842105263
$anchorCountNum.Copy($ws.Cells.Item(30,4))
$ws.Cells.Item(30,4).Value = 1
$anchorPctNum.Copy($ws.Cells.Item(30,5))
$ws.Cells.Item(30,5).Value = -100
$ws.Cells.Item(30,7).Value = 2
$ws.Cells.Item(30,10).Value = 2
$ws.Cells.Item(30,11).Value = -50
$ws.Cells.Item(30,14).Value = -94.117647058823
